$d = $word.ActiveDocument

function Insert-Xml($rng, [string]$bodyXml) {
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           $bodyXml + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
    $rng.InsertXML($pkg)
}

# ---------------------------------------------------------------------------
# 1) "...in the last fourth-quarter." -> split trailing words into own runs
# ---------------------------------------------------------------------------
$needle1 = "online sales revenue peaked in every fourth quarter of the year, which is a remarkable phenomenon worth to have further discussion. In first fourth-quarter, 43.7 and 8 billion U.S. dollars of sales revenue has gained by Amazon and Alibaba, and the value has been grown to 87.4 and 24.2 billion U.S. dollar in the last fourth-quarter."
$r1 = $d.Content
$ok1 = $r1.Find.Execute($needle1, $true)
if (-not $ok1) { throw "could not find paragraph about peaked sales revenue" }
$t1 = $d.Range($r1.Start, $r1.End)
$body1 = '<w:body><w:p>' +
    '<w:r><w:t xml:space="preserve">online sales revenue peaked in every fourth quarter of the year, which is a remarkable phenomenon worth to have further discussion. In first fourth-quarter, 43.7 and 8 billion U.S. dollars of sales revenue has gained by Amazon and Alibaba, and the value has been grown to 87.4 and 24.2 billion U.S. dollar in the last </w:t></w:r>' +
    '<w:r><w:t>fourth quarter</w:t></w:r>' +
    '<w:r><w:t>.</w:t></w:r>' +
    '</w:p></w:body>'
Insert-Xml $t1 $body1

# ---------------------------------------------------------------------------
# 2) "As shown ..." paragraph: drop its stray pPr, split the "improved..."
#    run, then add a brand-new paragraph discussing holiday consumption
#    (keeping the following "E-payment" bullet paragraph intact).
# ---------------------------------------------------------------------------
$r2start = $d.Content
$ok2a = $r2start.Find.Execute("As shown and mentioned", $true)
if (-not $ok2a) { throw "could not find 'As shown and mentioned' paragraph" }
$startPos = $r2start.Start

$r2end = $d.Content
$ok2b = $r2end.Find.Execute("E-payment", $true)
if (-not $ok2b) { throw "could not find 'E-payment' paragraph" }
$endPos = $r2end.Start + 1

$t2 = $d.Range($startPos, $endPos)

$rsquo = [char]0x2019
$ldquo = [char]0x201C
$rdquo = [char]0x201D

$p2a = '<w:p>' +
    '<w:r><w:t xml:space="preserve">As shown and mentioned in the graph and above, the reason of online sales increase is worth for discussion. The technology of online payment, as known as e-payment, has been </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">improved </w:t></w:r>' +
    '<w:r><w:t xml:space="preserve">and popularized </w:t></w:r>' +
    '<w:r><w:t>for years, which provides a perfect platform for online commercial activities.</w:t></w:r>' +
    '</w:p>'

$p2b = '<w:p>' +
    '<w:r><w:t>First, in every fourth quarter of the year, both Alibaba and Amazon' + $rsquo + 's online sales revenue increased</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to the top of the year. This phenomenon is usually come with holiday consumption, </w:t></w:r>' +
    '<w:r><w:rPr><w:color w:val="FF0000"/></w:rPr><w:t>explain the meaning of the word ' + $ldquo + 'holiday consumption' + $rdquo + '</w:t></w:r>' +
    '<w:r><w:t>. Significant festival such as Christmas</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> and the eleventh of November are highly triggers c</w:t></w:r>' +
    '<w:r><w:t>onsumer</w:t></w:r>' +
    '<w:r><w:t>' + $rsquo + 's</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> desire</w:t></w:r>' +
    '<w:r><w:t xml:space="preserve"> to buy products to treat themselves.</w:t></w:r>' +
    '</w:p>'

$p2c = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr>' +
    '<w:r><w:t>E-payment</w:t></w:r></w:p>'

$body2 = '<w:body>' + $p2a + $p2b + $p2c + '</w:body>'
Insert-Xml $t2 $body2

# ---------------------------------------------------------------------------
# 3) Move the lastRenderedPageBreak marker: it used to sit on "Conclusion ",
#    now the new content pushes the page break up to "E-payment" instead.
# ---------------------------------------------------------------------------
$r3c = $d.Content
$ok3c = $r3c.Find.Execute("Conclusion ", $true)
if (-not $ok3c) { throw "could not find 'Conclusion ' run" }
$t3c = $d.Range($r3c.Start, $r3c.End)
Insert-Xml $t3c '<w:body><w:p><w:r><w:t xml:space="preserve">Conclusion </w:t></w:r></w:p></w:body>'

$r3e = $d.Content
$ok3e = $r3e.Find.Execute("E-payment", $true)
if (-not $ok3e) { throw "could not find 'E-payment' run" }
$t3e = $d.Range($r3e.Start, $r3e.End)
Insert-Xml $t3e '<w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>E-payment</w:t></w:r></w:p></w:body>'

Write-Output "done"
